# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Estado de Cuenta" detail rows (16-27) are re-grouped/re-sorted by
# worker: KAREN SILGADO AYALA first, then all of EMIRO RAFAEL MARTINEZ
# BENITEZ's overdue periods in ascending order (2003..2008), then LUIS
# ANGEL ACOSTA HERNANDEZ, and finally the original first four workers
# (WILMER, IRAYDA, OSNAIDER, ERIS) moved to the bottom of the list.
# The underlying data set itself is unchanged - only the row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 16
$endRow = 27
$n = $endRow - $startRow + 1

# Read the current (pre-edit) contents of the data columns (C..G) for the
# 12 detail rows so we can rewrite them in the new order.
$colC = @()
$colD = @()
$colE = @()
$colF = @()
$colG = @()

for ($i = 0; $i -lt $n; $i++) {
    $r = $startRow + $i
    $colC += ,$ws.Cells.Item($r, 3).Value2
    $colD += ,$ws.Cells.Item($r, 4).Value2
    $colE += ,$ws.Cells.Item($r, 5).Value2
    $colF += ,$ws.Cells.Item($r, 6).Value2
    $colG += ,$ws.Cells.Item($r, 7).Value2
}

# new row i (0-based) = old row perm[i] (0-based), i.e. the row that used
# to sit at offset perm[i] below $startRow now sits at offset i.
$perm = @(4, 10, 9, 8, 7, 6, 5, 11, 0, 1, 2, 3)

for ($i = 0; $i -lt $n; $i++) {
    $r = $startRow + $i
    $p = $perm[$i]
    $ws.Cells.Item($r, 3).Value = $colC[$p]
    $ws.Cells.Item($r, 4).Value = $colD[$p]
    $ws.Cells.Item($r, 5).Value = $colE[$p]
    $ws.Cells.Item($r, 6).Value = $colF[$p]
    $ws.Cells.Item($r, 7).Value = $colG[$p]
}

"Reordered rows $startRow to $endRow"
